$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17+ down.
$ws.Rows.Item(17).Insert()

# The inserted row inherits formatting from the row above on columns B:F;
# that's not wanted here since only A and G are populated for the new entry.
$ws.Range("B17:F17").Clear()

# Fill in the new row's data: course name + rating.
$ws.Range("A17").Value = "Introduction to Excel"
$ws.Range("G17").Font.ThemeColor = 1
$ws.Range("G17").Value = 4
